$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the G column "chosenIntrinsicWidth" values that were hard-coded to
#    380 for a handful of rows but should reflect the actual chosen width
#    (matching column F for most of them).
# ---------------------------------------------------------------------------
$ws.Range("G9").Value  = 450
$ws.Range("G11").Value = 678
$ws.Range("G14").Value = 704
$ws.Range("G16").Value = 450

# ---------------------------------------------------------------------------
# 2. Tweak the evaluation formula in column J:
#      "POOR! (--)" -> "POOR (--)"
#      "OK"         -> "GOOD!"
# ---------------------------------------------------------------------------
$ws.Range("J2:J17").Formula = '=IFS(I2<0.9, "POOR (--)", I2<1, "(-)", I2=1, "GOOD!", I2>1.2, "BIG (++)", I2>1, "(+)")'

# ---------------------------------------------------------------------------
# 3. Remove the static yellow fill that used to highlight column G - it is
#    being replaced by the conditional formatting added below.
# ---------------------------------------------------------------------------
$ws.Range("G1:G17").Interior.Pattern = -4142   # xlNone

# ---------------------------------------------------------------------------
# 4. Add conditional formatting to G2:G17 driven by the evaluation in
#    column J: green when "GOOD!", red shades for "POOR (--)"/"BIG (++)",
#    light green for the milder "(+)"/"(-)" results.
#    Rules are added in priority order (1-4) first, then each one's dxf
#    (fill/font) is assigned starting from the last rule backwards so the
#    resulting dxfId allocation matches what Excel produces when rules are
#    authored through the "Conditional Formatting" dialog.
# ---------------------------------------------------------------------------
$rng = $ws.Range("G2:G17")

$fcGood = $rng.FormatConditions.Add(2, 3, 'J2="GOOD!"')
$fcPoor = $rng.FormatConditions.Add(2, 3, 'J2="POOR (--)"')
$fcBig  = $rng.FormatConditions.Add(2, 3, 'J2="BIG (++)"')
$fcMild = $rng.FormatConditions.Add(2, 3, 'OR(J2="(+)",J2="(-)")')

$fcMild.Interior.Color = 13823444   # FFD4EDD2 light green

$fcBig.Font.Color = 16777215        # white
$fcBig.Interior.Color = 192         # FFC00000 dark red

$fcPoor.Font.Color = 16777215       # white
$fcPoor.Interior.Color = 5066944    # FFC0504D medium red

$fcGood.Interior.Color = 8048793    # FF99D07A green

# ---------------------------------------------------------------------------
# 5. Restore the A1:K1 selection state.
# ---------------------------------------------------------------------------
[void]$ws.Range("A1:K1").Select()
